$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.142.45"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "2.307.35"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.54%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.521"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.44%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.517"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.85"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0790"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.38%  "
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.04"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.89"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.41%  "
$ws.Range("D15").Value = "2.666.21"
$ws.Range("E15").Value = "  -0.04%  "
$ws.Range("D16").Value = "2.254.32"
$ws.Range("E16").Value = "  -1.31%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.789"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.08%  "
$ws.Range("D18").Value = "43.044.69"
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.31"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.40%  "
$ws.Range("D20").Value = "0.0₃0909"
$ws.Range("E20").Value = "  +0.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.30"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.55%  "
$ws.Range("E24").Value = "  -2.58%  "
$ws.Range("E25").Value = "  -0.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.43"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.84"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "168.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.95%  "
$ws.Range("E29").Value = "  -0.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.03"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.91"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.33%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.07%  "
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.84"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.10"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.03%  "
$ws.Range("E36").Value = "  -0.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0689"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.23%  "
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.79"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.64%  "
$ws.Range("E40").Value = "  +1.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.76"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.66%  "
$ws.Range("D42").Value = "2.004.96"
$ws.Range("E42").Value = "  +0.42%  "
$ws.Range("E43").Value = "  -0.92%  "
$ws.Range("E44").Value = "  -4.55%  "
$ws.Range("E45").Value = "  -1.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.31"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.83"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.42"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.29%  "
$ws.Range("D49").Value = "2.528.45"
$ws.Range("E49").Value = "  -0.25%  "
$ws.Range("E50").Value = "  -1.67%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.63"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.51%  "
